$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 326.42856
$ws.Range("J2").Value = 745
$ws.Range("L2").Value = 745
$ws.Range("N2").Value = -971
$ws.Range("H62").Value = 3644.9
$ws.Range("I62").Value = 3272.111
$ws.Range("K62").Value = 3272.111
$ws.Range("M62").Value = -2648.111
$ws.Range("H65").Value = 3644.9
$ws.Range("I65").Value = 3272.111
$ws.Range("K65").Value = 16360.555
$ws.Range("M65").Value = -13240.555
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").Value = $null
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").Value = $null
$ws.Range("H138").Value = 3686.0205
$ws.Range("I138").Value = 3327.6667
$ws.Range("J138").Value = 3766.65
$ws.Range("K138").Value = 9983.000100000001
$ws.Range("L138").Value = 11299.95
$ws.Range("M138").Value = -4843.000100000001
$ws.Range("N138").Value = -21579.95
$ws.Range("H141").Value = 4415.3125
$ws.Range("I141").Value = 4431.8184
$ws.Range("K141").Value = 13295.4552
$ws.Range("M141").Value = -8115.4552

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2443199.8
$ws.Range("I32").Value = 2780280
$ws.Range("K32").Value = 2780280
$ws.Range("M32").Value = -2779993
$ws.Range("H88").Value = 1566.7894
$ws.Range("I88").Value = 821.6667
$ws.Range("J88").Value = 1706.5
$ws.Range("K88").Value = 821.6667
$ws.Range("L88").Value = 1706.5
$ws.Range("M88").Value = -415.6667
$ws.Range("N88").Value = -2518.5
$ws.Range("H91").Value = 1566.7894
$ws.Range("I91").Value = 821.6667
$ws.Range("J91").Value = 1706.5
$ws.Range("K91").Value = 821.6667
$ws.Range("L91").Value = 1706.5
$ws.Range("M91").Value = 582.3333
$ws.Range("N91").Value = -4514.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H112").Value = 103000
$ws.Range("J112").Value = 103000
$ws.Range("L112").Value = 103000
$ws.Range("N112").Value = -105954

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 30756
$ws.Range("I94").Value = 60012
$ws.Range("J94").Value = 1500
$ws.Range("K94").Value = 60012
$ws.Range("L94").Value = 1500
$ws.Range("M94").Value = -59561
$ws.Range("N94").Value = -2402

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 37.875
$ws.Range("I33").Value = 10.75
$ws.Range("J33").Value = 65
$ws.Range("K33").Value = 64.5
$ws.Range("L33").Value = 390
$ws.Range("M33").Value = 218.5
$ws.Range("N33").Value = -956
$ws.Range("H113").Value = 3715.9443
$ws.Range("I113").Value = 4661
$ws.Range("K113").Value = 13983
$ws.Range("M113").Value = -11813
$ws.Range("H131").Value = 1516.41
$ws.Range("I131").Value = 1266
$ws.Range("J131").Value = 1529.5895
$ws.Range("K131").Value = 3798
$ws.Range("L131").Value = 4588.7685
$ws.Range("M131").Value = 1242
$ws.Range("N131").Value = -14668.7685

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 968
$ws.Range("I97").Value = 920.8421
$ws.Range("K97").Value = 920.8421
$ws.Range("M97").Value = -424.8421
$ws.Range("H122").Value = 3325.625
$ws.Range("I122").Value = 3325.625
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 9976.875
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -7526.875
$ws.Range("N122").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1661.4286
$ws.Range("I22").Value = 852.2727
$ws.Range("J22").Value = 2551.5
$ws.Range("K22").Value = 852.2727
$ws.Range("L22").Value = 2551.5
$ws.Range("M22").Value = -557.2727
$ws.Range("N22").Value = -3141.5
$ws.Range("H25").Value = 10000
$ws.Range("I25").Value = 10000
$ws.Range("K25").Value = 10000
$ws.Range("M25").Value = -9770
$ws.Range("H27").Value = 1661.4286
$ws.Range("I27").Value = 852.2727
$ws.Range("J27").Value = 2551.5
$ws.Range("K27").Value = 852.2727
$ws.Range("L27").Value = 2551.5
$ws.Range("M27").Value = -745.2727
$ws.Range("N27").Value = -2765.5
$ws.Range("H38").Value = 69763.46000000001
$ws.Range("I38").Value = 48015.668
$ws.Range("J38").Value = 76287.8
$ws.Range("K38").Value = 48015.668
$ws.Range("L38").Value = 76287.8
$ws.Range("M38").Value = -47605.668
$ws.Range("N38").Value = -77107.8
$ws.Range("H40").Value = 4996
$ws.Range("I40").Value = 4980
$ws.Range("K40").Value = 4980
$ws.Range("M40").Value = -4844
$ws.Range("H46").Value = 2215.8333
$ws.Range("I46").Value = 1850
$ws.Range("J46").Value = 2398.75
$ws.Range("K46").Value = 1850
$ws.Range("L46").Value = 2398.75
$ws.Range("M46").Value = -1662
$ws.Range("N46").Value = -2774.75
$ws.Range("H76").Value = 63821.25
$ws.Range("J76").Value = 63821.25
$ws.Range("L76").Value = 63821.25
$ws.Range("N76").Value = -64497.25
$ws.Range("H79").Value = 63821.25
$ws.Range("J79").Value = 63821.25
$ws.Range("L79").Value = 63821.25
$ws.Range("N79").Value = -66161.25
$ws.Range("H82").Value = 1863.8387
$ws.Range("I82").Value = 1989.5
$ws.Range("J82").Value = 1635.3636
$ws.Range("K82").Value = 1989.5
$ws.Range("L82").Value = 1635.3636
$ws.Range("M82").Value = -1628.5
$ws.Range("N82").Value = -2357.3636
$ws.Range("H85").Value = 1863.8387
$ws.Range("I85").Value = 1989.5
$ws.Range("J85").Value = 1635.3636
$ws.Range("K85").Value = 1989.5
$ws.Range("L85").Value = 1635.3636
$ws.Range("M85").Value = -741.5
$ws.Range("N85").Value = -4131.3636
$ws.Range("H93").Value = 3762.2778
$ws.Range("J93").Value = 738.6
$ws.Range("L93").Value = 738.6
$ws.Range("N93").Value = -3234.6
$ws.Range("H94").Value = 29330
$ws.Range("J94").Value = 29330
$ws.Range("L94").Value = 29330
$ws.Range("N94").Value = -30682

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4057.8
$ws.Range("I81").Value = 4064.4443
$ws.Range("K81").Value = 8128.8886
$ws.Range("M81").Value = -7067.8886
$ws.Range("H84").Value = 4057.8
$ws.Range("I84").Value = 4064.4443
$ws.Range("K84").Value = 40644.443
$ws.Range("M84").Value = -35340.443
$ws.Range("H111").Value = 49000
$ws.Range("J111").Value = 49000
$ws.Range("L111").Value = 49000
$ws.Range("N111").Value = -57180
$ws.Range("H122").Value = 2399
$ws.Range("I122").Value = 1959.8
$ws.Range("J122").Value = 3497
$ws.Range("K122").Value = 5879.4
$ws.Range("L122").Value = 10491
$ws.Range("M122").Value = -3429.4
$ws.Range("N122").Value = -15391
$ws.Range("H126").Value = 3198.158
$ws.Range("I126").Value = 3364.1177
$ws.Range("J126").Value = 1787.5
$ws.Range("K126").Value = 10092.3531
$ws.Range("L126").Value = 5362.5
$ws.Range("M126").Value = -7622.3531
$ws.Range("N126").Value = -10302.5
$ws.Range("H136").Value = 336736.8
$ws.Range("J136").Value = 1486802
$ws.Range("L136").Value = 4460406
$ws.Range("N136").Value = -4465506
